# Daily attendance processing - 2026-01-10 20:35:58
# Reorders the "Recorded By" column (G) entries so that "System" is listed
# first for the known automated-recorder identities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq "backup@backdoor.com, System") {
        $cell.Value = "System, backup@backdoor.com"
    }
    elseif ($value -eq "backup@backdoor.com, system, System") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
    elseif ($value -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
}
